# Insert two new weekly price-report rows (Pomelo, Start Ruby) for the
# Feria Lagunitas de Puerto Montt market, dated 44769 (2022-07-27), right
# before the existing row 309. This pushes the old rows 309:339 down to
# 311:341 (Excel handles that shift automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("309:310").Insert()

# --- New row 309 (Primera) ---
$ws.Range("A309").Value = 4
$ws.Range("B309").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C309").Value = "Los Lagos"
$ws.Range("D309").Value = 44769
$ws.Range("E309").Value = 10
$ws.Range("F309").Value = "Fruta"
$ws.Range("G309").Value = 100102
$ws.Range("H309").Value = "Cítricos"
$ws.Range("I309").Value = 100102006
$ws.Range("J309").Value = "Pomelo"
$ws.Range("K309").Value = "Start Ruby"
$ws.Range("L309").Value = "Primera"
$ws.Range("M309").Value = 40
$ws.Range("N309").Value = 14000
$ws.Range("O309").Value = 15000
$ws.Range("P309").Value = 14500
$ws.Range("Q309").Value = "$/caja 14 kilos empedrada"
$ws.Range("R309").Value = "Región de O'Higgins"
$ws.Range("S309").Value = 1036
$ws.Range("T309").Value = 14

# --- New row 310 (Segunda) ---
$ws.Range("A310").Value = 4
$ws.Range("B310").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C310").Value = "Los Lagos"
$ws.Range("D310").Value = 44769
$ws.Range("E310").Value = 10
$ws.Range("F310").Value = "Fruta"
$ws.Range("G310").Value = 100102
$ws.Range("H310").Value = "Cítricos"
$ws.Range("I310").Value = 100102006
$ws.Range("J310").Value = "Pomelo"
$ws.Range("K310").Value = "Start Ruby"
$ws.Range("L310").Value = "Segunda"
$ws.Range("M310").Value = 20
$ws.Range("N310").Value = 12000
$ws.Range("O310").Value = 12000
$ws.Range("P310").Value = 12000
$ws.Range("Q310").Value = "$/caja 14 kilos empedrada"
$ws.Range("R310").Value = "Región de O'Higgins"
$ws.Range("S310").Value = 857
$ws.Range("T310").Value = 14
